$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1442.9
$ws.Range("J40").Value = 1725.5714
$ws.Range("L40").Value = 1725.5714
$ws.Range("N40").Value = -2075.5714
$ws.Range("H53").Value = 2348.5715
$ws.Range("I53").Value = 347
$ws.Range("J53").Value = 2894.4546
$ws.Range("K53").Value = 347
$ws.Range("L53").Value = 2894.4546
$ws.Range("M53").Value = 290
$ws.Range("N53").Value = -4168.4546
$ws.Range("H70").Value = 1560
$ws.Range("I70").Value = 1433.3334
$ws.Range("K70").Value = 4300.0002
$ws.Range("M70").Value = -4030.0002
$ws.Range("H73").Value = 1560
$ws.Range("I73").Value = 1433.3334
$ws.Range("K73").Value = 4300.0002
$ws.Range("M73").Value = -3364.0002
$ws.Range("H88").Value = 398
$ws.Range("I88").Value = 122.5
$ws.Range("J88").Value = 1500
$ws.Range("K88").Value = 122.5
$ws.Range("L88").Value = 1500
$ws.Range("M88").Value = 283.5
$ws.Range("N88").Value = -2312
$ws.Range("H91").Value = 398
$ws.Range("I91").Value = 122.5
$ws.Range("J91").Value = 1500
$ws.Range("K91").Value = 122.5
$ws.Range("L91").Value = 1500
$ws.Range("M91").Value = 1281.5
$ws.Range("N91").Value = -4308
$ws.Range("H94").Value = 2970
$ws.Range("I94").Value = 2970
$ws.Range("K94").Value = 2970
$ws.Range("M94").Value = -2519
$ws.Range("H103").Value = 250187.8
$ws.Range("I103").Value = 500244.4
$ws.Range("J103").Value = 131.2
$ws.Range("K103").Value = 1500733.2
$ws.Range("L103").Value = 393.6
$ws.Range("M103").Value = -1500147.2
$ws.Range("N103").Value = -1565.6
$ws.Range("H106").Value = 7938707.5
$ws.Range("I106").Value = 11495801
$ws.Range("K106").Value = 11495801
$ws.Range("M106").Value = -11495170
$ws.Range("H116").Value = 8000.857
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 8667.666999999999
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 8667.666999999999
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -15551.667
$ws.Range("H125").Value = 1108.75
$ws.Range("J125").Value = 1345
$ws.Range("L125").Value = 12105
$ws.Range("N125").Value = -17025
$ws.Range("H135").Value = 17859222
$ws.Range("I135").Value = 572.65216
$ws.Range("K135").Value = 5153.869439999999
$ws.Range("M135").Value = -2618.869439999999
$ws.Range("H138").Value = 30305906
$ws.Range("I138").Value = 76925016
$ws.Range("J138").Value = 3484.15
$ws.Range("K138").Value = 230775048
$ws.Range("L138").Value = 10452.45
$ws.Range("M138").Value = -230769908
$ws.Range("N138").Value = -20732.45

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 500
$ws.Range("J10").Value = 500
$ws.Range("L10").Value = 500
$ws.Range("N10").Value = -840
$ws.Range("H61").Value = 4028.6843
$ws.Range("I61").Value = 4002.647
$ws.Range("K61").Value = 4002.647
$ws.Range("M61").Value = -3790.647
$ws.Range("H74").Value = 45456716
$ws.Range("I74").Value = 100001730
$ws.Range("J74").Value = 2542.8333
$ws.Range("K74").Value = 100001730
$ws.Range("L74").Value = 2542.8333
$ws.Range("M74").Value = -100000856
$ws.Range("N74").Value = -4290.8333
$ws.Range("H77").Value = 45456716
$ws.Range("I77").Value = 100001730
$ws.Range("J77").Value = 2542.8333
$ws.Range("K77").Value = 500008650
$ws.Range("L77").Value = 12714.1665
$ws.Range("M77").Value = -500004282
$ws.Range("N77").Value = -21450.1665
$ws.Range("H132").Value = 34763.5
$ws.Range("I132").Value = 3547.111
$ws.Range("J132").Value = 74898.86
$ws.Range("K132").Value = 10641.333
$ws.Range("L132").Value = 224696.58
$ws.Range("M132").Value = -8111.332999999999
$ws.Range("N132").Value = -229756.58
$ws.Range("H136").Value = 4028.6843
$ws.Range("I136").Value = 4002.647
$ws.Range("K136").Value = 12007.941
$ws.Range("M136").Value = -9457.940999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1523.439
$ws.Range("I86").Value = 1438.44
$ws.Range("J86").Value = 1656.25
$ws.Range("K86").Value = 1438.44
$ws.Range("L86").Value = 1656.25
$ws.Range("M86").Value = -315.4400000000001
$ws.Range("N86").Value = -3902.25
$ws.Range("H89").Value = 1523.439
$ws.Range("I89").Value = 1438.44
$ws.Range("J89").Value = 1656.25
$ws.Range("K89").Value = 7192.200000000001
$ws.Range("L89").Value = 8281.25
$ws.Range("M89").Value = -1576.200000000001
$ws.Range("N89").Value = -19513.25
$ws.Range("H105").Value = 7145628.5
$ws.Range("I105").Value = 2550
$ws.Range("K105").Value = 2550
$ws.Range("M105").Value = -803
$ws.Range("H110").Value = 20702
$ws.Range("J110").Value = 20702
$ws.Range("L110").Value = 20702
$ws.Range("N110").Value = -28882
$ws.Range("H134").Value = 3836.9285
$ws.Range("I134").Value = 3942
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 11826
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -9291
$ws.Range("N134").Value = -8070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 767.2
$ws.Range("J12").Value = 1125.3334
$ws.Range("L12").Value = 1125.3334
$ws.Range("N12").Value = -1465.3334
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736
$ws.Range("H96").Value = 12624
$ws.Range("J96").Value = 12624
$ws.Range("L96").Value = 12624
$ws.Range("N96").Value = -18116
$ws.Range("H122").Value = 2345
$ws.Range("I122").Value = 2345
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7035
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -4585
$ws.Range("H132").Value = 2860.889
$ws.Range("I132").Value = 1202.6154
$ws.Range("J132").Value = 7172.4
$ws.Range("K132").Value = 3607.8462
$ws.Range("L132").Value = 21517.2
$ws.Range("M132").Value = -1077.8462
$ws.Range("N132").Value = -26577.2
$ws.Range("H134").Value = 1398.9445
$ws.Range("I134").Value = 1339.4286
$ws.Range("J134").Value = 1607.25
$ws.Range("K134").Value = 4018.2858
$ws.Range("L134").Value = 4821.75
$ws.Range("M134").Value = -1483.2858
$ws.Range("N134").Value = -9891.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 271
$ws.Range("I108").Value = 271
$ws.Range("K108").Value = 813
$ws.Range("M108").Value = 2067
$ws.Range("H131").Value = 715.12
$ws.Range("J131").Value = 723.5625
$ws.Range("L131").Value = 2170.6875
$ws.Range("N131").Value = -12250.6875
$ws.Range("H141").Value = 1866.25
$ws.Range("I141").Value = 1866.25
$ws.Range("K141").Value = 5598.75
$ws.Range("M141").Value = -418.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3192.2104
$ws.Range("I102").Value = 3316
$ws.Range("J102").Value = 2845.6
$ws.Range("K102").Value = 3316
$ws.Range("L102").Value = 2845.6
$ws.Range("M102").Value = -1694
$ws.Range("N102").Value = -6089.6
$ws.Range("H122").Value = 5666.6665
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550
$ws.Range("H126").Value = 5264.077
$ws.Range("I126").Value = 4379.4287
$ws.Range("J126").Value = 6296.1665
$ws.Range("K126").Value = 13138.2861
$ws.Range("L126").Value = 18888.4995
$ws.Range("M126").Value = -10668.2861
$ws.Range("N126").Value = -23828.4995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9150.909
$ws.Range("I7").Value = 4870
$ws.Range("J7").Value = 20566.666
$ws.Range("K7").Value = 4870
$ws.Range("L7").Value = 20566.666
$ws.Range("M7").Value = -4758
$ws.Range("N7").Value = -20790.666
$ws.Range("H68").Value = 2127.1667
$ws.Range("I68").Value = 1550
$ws.Range("J68").Value = 2415.75
$ws.Range("K68").Value = 1550
$ws.Range("L68").Value = 2415.75
$ws.Range("M68").Value = -801
$ws.Range("N68").Value = -3913.75
$ws.Range("H71").Value = 2127.1667
$ws.Range("I71").Value = 1550
$ws.Range("J71").Value = 2415.75
$ws.Range("K71").Value = 7750
$ws.Range("L71").Value = 12078.75
$ws.Range("M71").Value = -4006
$ws.Range("N71").Value = -19566.75
$ws.Range("H82").Value = 3928.5715
$ws.Range("I82").Value = 4166.6665
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 4166.6665
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -3805.6665
$ws.Range("N82").Value = -3222
$ws.Range("H85").Value = 3928.5715
$ws.Range("I85").Value = 4166.6665
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 4166.6665
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -2918.6665
$ws.Range("N85").Value = -4996
$ws.Range("H122").Value = 1403702.8
$ws.Range("I122").Value = 2181537.5
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 6544612.5
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -6542162.5
$ws.Range("N122").Value = -15700
$ws.Range("H126").Value = 9150.909
$ws.Range("I126").Value = 4870
$ws.Range("J126").Value = 20566.666
$ws.Range("K126").Value = 14610
$ws.Range("L126").Value = 61699.99800000001
$ws.Range("M126").Value = -12140
$ws.Range("N126").Value = -66639.99800000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4450.25
$ws.Range("I62").Value = 3200.6667
$ws.Range("J62").Value = 5200
$ws.Range("K62").Value = 3200.6667
$ws.Range("L62").Value = 5200
$ws.Range("M62").Value = -2576.6667
$ws.Range("N62").Value = -6448
$ws.Range("H65").Value = 4450.25
$ws.Range("I65").Value = 3200.6667
$ws.Range("J65").Value = 5200
$ws.Range("K65").Value = 16003.3335
$ws.Range("L65").Value = 26000
$ws.Range("M65").Value = -12883.3335
$ws.Range("N65").Value = -32240
$ws.Range("H132").Value = 1192.1724
$ws.Range("I132").Value = 726.4737
$ws.Range("J132").Value = 2077
$ws.Range("K132").Value = 2179.4211
$ws.Range("L132").Value = 6231
$ws.Range("M132").Value = 350.5789
$ws.Range("N132").Value = -11291
$ws.Range("H136").Value = 31253340
$ws.Range("I136").Value = 41668160
$ws.Range("J136").Value = 8888.125
$ws.Range("K136").Value = 125004480
$ws.Range("L136").Value = 26664.375
$ws.Range("M136").Value = -125001930
$ws.Range("N136").Value = -31764.375
$ws.Range("H140").Value = 39025.8
$ws.Range("J140").Value = 39025.8
$ws.Range("L140").Value = 39025.8
$ws.Range("N140").Value = -49385.8
